# Bulk Upload Template — complete bulk upload sample data for role and
# technology sheets, and restore the Technologies tab as the active sheet.

$wb = $excel.ActiveWorkbook

$wsTech  = $wb.Worksheets.Item("Technologies")
$wsEmp   = $wb.Worksheets.Item("Employees")
$wsRoles = $wb.Worksheets.Item("Roles")

# ---------------------------------------------------------------------
# Technologies sheet: drop the ".Net" sample row and replace the
# trailing "Java" sample with another "CSS" entry.
# ---------------------------------------------------------------------
$wsTech.Range("B3").ClearContents() | Out-Null
$wsTech.Range("B7").Value = "CSS"

# ---------------------------------------------------------------------
# Employees sheet: no data changes — just move the selection.
# ---------------------------------------------------------------------
$wsEmp.Activate() | Out-Null
$wsEmp.Range("B9").Select() | Out-Null

# ---------------------------------------------------------------------
# Roles sheet: drop "Data engineer 2" and "Data engineer 6" rows, and
# replace "Data engineer 5" with a duplicate "Data engineer 4".
# ---------------------------------------------------------------------
$wsRoles.Activate() | Out-Null
$wsRoles.Range("B3").ClearContents() | Out-Null
$wsRoles.Range("B5").Value = "Data engineer 4"
$wsRoles.Range("B6").ClearContents() | Out-Null
$wsRoles.Range("H14").Select() | Out-Null

# ---------------------------------------------------------------------
# Technologies becomes the active tab again (matches activeTab="0").
# ---------------------------------------------------------------------
$wsTech.Activate() | Out-Null
$wsTech.Range("B7").Select() | Out-Null
